$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1200
$ws.Range("J40").Value = 1200
$ws.Range("L40").Value = 1200
$ws.Range("N40").Value = -1550
$ws.Range("H62").Value = 2000
$ws.Range("J62").Value = 2000
$ws.Range("L62").Value = 2000
$ws.Range("N62").Value = -3248
$ws.Range("H65").Value = 2000
$ws.Range("J65").Value = 2000
$ws.Range("L65").Value = 10000
$ws.Range("N65").Value = -16240
$ws.Range("H100").Value = 1309
$ws.Range("I100").Value = 1286
$ws.Range("J100").Value = 1370.3334
$ws.Range("K100").Value = 1286
$ws.Range("L100").Value = 1370.3334
$ws.Range("M100").Value = -745
$ws.Range("N100").Value = -2452.3334
$ws.Range("H135").Value = 887
$ws.Range("I135").Value = 882.6667
$ws.Range("K135").Value = 7944.0003
$ws.Range("M135").Value = -5409.0003
$ws.Range("H141").Value = 2498.25
$ws.Range("I141").Value = 2498.25
$ws.Range("K141").Value = 7494.75
$ws.Range("M141").Value = -2314.75

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1437.6316
$ws.Range("I2").Value = 901.05884
$ws.Range("J2").Value = 5998.5
$ws.Range("K2").Value = 901.05884
$ws.Range("L2").Value = 5998.5
$ws.Range("M2").Value = -788.05884
$ws.Range("N2").Value = -6224.5
$ws.Range("H6").Value = 1534
$ws.Range("I6").Value = 1534
$ws.Range("K6").Value = 1534
$ws.Range("M6").Value = -1361
$ws.Range("H61").Value = 1492
$ws.Range("I61").Value = 1492
$ws.Range("K61").Value = 1492
$ws.Range("M61").Value = -1280
$ws.Range("H116").Value = 1437.6316
$ws.Range("I116").Value = 901.05884
$ws.Range("J116").Value = 5998.5
$ws.Range("K116").Value = 901.05884
$ws.Range("L116").Value = 5998.5
$ws.Range("M116").Value = 1392.94116
$ws.Range("N116").Value = -10586.5
$ws.Range("H122").Value = 1798.5
$ws.Range("I122").Value = 1798.5
$ws.Range("K122").Value = 5395.5
$ws.Range("M122").Value = -2945.5
$ws.Range("H132").Value = 839.4
$ws.Range("I132").Value = 858.6923
$ws.Range("J132").Value = 714
$ws.Range("K132").Value = 2576.0769
$ws.Range("L132").Value = 2142
$ws.Range("M132").Value = -46.07690000000002
$ws.Range("N132").Value = -7202
$ws.Range("H136").Value = 1492
$ws.Range("I136").Value = 1492
$ws.Range("K136").Value = 4476
$ws.Range("M136").Value = -1926

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1437.6316
$ws.Range("I3").Value = 901.05884
$ws.Range("J3").Value = 5998.5
$ws.Range("K3").Value = 901.05884
$ws.Range("L3").Value = 5998.5
$ws.Range("M3").Value = -787.05884
$ws.Range("N3").Value = -6226.5
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 100
$ws.Range("K22").Value = 100
$ws.Range("M22").Value = 73
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9011.5
$ws.Range("J31").Value = 9011
$ws.Range("L31").Value = 9011
$ws.Range("N31").Value = -9601
$ws.Range("H34").Value = 9011.5
$ws.Range("J34").Value = 9011
$ws.Range("L34").Value = 9011
$ws.Range("N34").Value = -9415
$ws.Range("H42").Value = 8000
$ws.Range("I42").Value = 8000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 8000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -7407
$ws.Range("N42").ClearContents()
$ws.Range("H107").Value = 710.2105
$ws.Range("J107").Value = 450
$ws.Range("L107").Value = 450
$ws.Range("N107").Value = -4290
$ws.Range("H132").Value = 1379.2858
$ws.Range("I132").Value = 1379.2858
$ws.Range("K132").Value = 4137.857400000001
$ws.Range("M132").Value = -1607.857400000001

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3664.4443
$ws.Range("J4").Value = 6300
$ws.Range("L4").Value = 18900
$ws.Range("N4").Value = -19124
$ws.Range("H11").Value = 4400
$ws.Range("J11").Value = 5249.75
$ws.Range("L11").Value = 15749.25
$ws.Range("N11").Value = -16029.25
$ws.Range("H13").Value = 10015.833
$ws.Range("I13").Value = 47.5
$ws.Range("J13").Value = 15000
$ws.Range("K13").Value = 142.5
$ws.Range("L13").Value = 45000
$ws.Range("M13").Value = 25.5
$ws.Range("N13").Value = -45336
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("M121").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 129.08333
$ws.Range("J2").Value = 152
$ws.Range("L2").Value = 152
$ws.Range("N2").Value = -378
$ws.Range("H5").Value = 12699.5
$ws.Range("I5").Value = 12699.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 12699.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -12587.5
$ws.Range("N5").ClearContents()
$ws.Range("H122").Value = 1544
$ws.Range("I122").Value = 1544
$ws.Range("K122").Value = 4632
$ws.Range("M122").Value = -2182

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 11000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H22").Value = 1910.5
$ws.Range("I22").Value = 880.6
$ws.Range("J22").Value = 2940.4
$ws.Range("K22").Value = 880.6
$ws.Range("L22").Value = 2940.4
$ws.Range("M22").Value = -585.6
$ws.Range("N22").Value = -3530.4
$ws.Range("H27").Value = 1910.5
$ws.Range("I27").Value = 880.6
$ws.Range("J27").Value = 2940.4
$ws.Range("K27").Value = 880.6
$ws.Range("L27").Value = 2940.4
$ws.Range("M27").Value = -773.6
$ws.Range("N27").Value = -3154.4
$ws.Range("H120").Value = 80465.336
$ws.Range("J120").Value = 80465.336
$ws.Range("L120").Value = 80465.336
$ws.Range("N120").Value = -90141.336
$ws.Range("H132").Value = 5322.6
$ws.Range("I132").Value = 1832.3334
$ws.Range("K132").Value = 5497.0002
$ws.Range("M132").Value = -2967.0002

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 867.2941
$ws.Range("I122").Value = 914.5
$ws.Range("K122").Value = 2743.5
$ws.Range("M122").Value = -293.5
$ws.Range("H125").Value = 34998
$ws.Range("J125").Value = 34998
$ws.Range("L125").Value = 34998
$ws.Range("N125").Value = -44838
$ws.Range("H136").Value = 948.75
$ws.Range("I136").Value = 963.3333
$ws.Range("J136").Value = 905
$ws.Range("K136").Value = 2889.9999
$ws.Range("L136").Value = 2715
$ws.Range("M136").Value = -339.9998999999998
$ws.Range("N136").Value = -7815
